# historico_vendas.xlsx edit
# - Remove the last two sale records (rows 3 and 4), shrinking the used range to A1:F2
# - Replace the remaining sale record (row 2) with a new entry

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 and 4 entirely (data + dimension shrink to A1:F2)
$ws.Rows("3:4").Delete()

# Update row 2 with the new sale record
$ws.Range("A2").Value = 45604.79457642078
$ws.Range("B2").Value = "Camiseta Estampa Animal"
$ws.Range("C2").Value = 1234
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 45
$ws.Range("F2").Value = 45
